$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 218; this pushes the former row 218
# (and everything below it) down by one, matching the target diff where
# former rows 218-281 become 219-282 and a brand-new row 218 is introduced.
$ws.Rows.Item(218).Insert()

# Populate the newly inserted row 218 with its data.
$ws.Range("A218").Value = 9
$ws.Range("B218").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C218").Value = "Metropolitana"
$ws.Range("D218").Value = 44736
$ws.Range("E218").Value = 13
$ws.Range("F218").Value = 300000001
$ws.Range("G218").Value = "Rabanito"
$ws.Range("H218").Value = "Sin especificar"
$ws.Range("I218").Value = "Primera"
$ws.Range("J218").Value = 7000
$ws.Range("K218").Value = 3000
$ws.Range("L218").Value = 3000
$ws.Range("M218").Value = 3000
$ws.Range("N218").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O218").Value = "Provincia de Chacabuco"
$ws.Range("P218").Value = 30
$ws.Range("Q218").Value = 100
$ws.Range("R218").Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of
# column D (style index 2 in the original workbook).
$ws.Range("D218").NumberFormat = $ws.Range("D219").NumberFormat
